$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column; existing columns A-D shift to B-E,
# carrying their exact widths/content/styles along with them.
$ws.Columns("A").Insert()

# New column A content (TabName / CasesTab)
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Narrow best-fit width for the new column A
$ws.Columns("A").ColumnWidth = 8.0

# Replace the two Cypher query cells in row 2 with their updated text
$casesQuery = 'MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.ethnicity IN [''NOT_HISPANIC_OR_LATINO''] 
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '''') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '''') AS `Trial Code`,
    COALESCE(a.arm_id, '''') AS `Arm`,
    COALESCE(a.arm_drug, '''') AS `Arm Treatment`,
    COALESCE(c.disease, '''') AS `Diagnosis`,
    COALESCE(c.gender, '''') AS `Gender`,
    COALESCE(c.race, '''') AS `Race`,
    COALESCE(c.ethnicity, '''') AS `Ethnicity`'
$statQuery = 'MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE WHERE c.ethnicity IN [''NOT_HISPANIC_OR_LATINO''] 
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials'
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# Row 2 grew taller to accommodate the longer wrapped query text
$ws.Rows(2).RowHeight = 174

# Selection moved to C2 in the saved workbook
$ws.Range("C2").Select() | Out-Null
